$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 0.74
$ws.Range("J2").Value = 0.6899999999999999
$ws.Range("M2").Value = 0.82
$ws.Range("P2").Value = 0.92

$ws.Range("D3").Value = 0.98
$ws.Range("G3").Value = 0.66
$ws.Range("J3").Value = 0.6
$ws.Range("M3").Value = 0.75
$ws.Range("P3").Value = 0.87

$ws.Range("D4").Value = 0.86
$ws.Range("G4").Value = 0.79
$ws.Range("J4").Value = 0.85
$ws.Range("M4").Value = 0.86
$ws.Range("P4").Value = 0.86

$ws.Range("D5").Value = 1
$ws.Range("G5").Value = 0.73
$ws.Range("J5").Value = 0.68
$ws.Range("M5").Value = 0.8100000000000001
$ws.Range("P5").Value = 0.91

$ws.Range("D6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("P6").ClearContents()

$ws.Range("D7").Value = 0
$ws.Range("G7").Value = 0.15
$ws.Range("J7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("P7").Value = 0

$ws.Range("D8").Value = 0
$ws.Range("G8").Value = 0.15
$ws.Range("J8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("P8").Value = 0

$ws.Range("D9").Value = 0
$ws.Range("G9").Value = 0.15
$ws.Range("J9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("P9").Value = 0

$ws.Range("D10").Value = 1
$ws.Range("G10").Value = 0.73
$ws.Range("J10").Value = 0.68
$ws.Range("M10").Value = 0.8100000000000001
$ws.Range("P10").Value = 0.91

$ws.Range("D11").Value = 0.97
$ws.Range("G11").Value = 0.77
$ws.Range("J11").Value = 0.74
$ws.Range("M11").Value = 0.84
$ws.Range("P11").Value = 0.91

$ws.Range("D12").Value = 0.92
$ws.Range("G12").Value = 0.39
$ws.Range("J12").Value = 0.29
$ws.Range("M12").Value = 0.44
$ws.Range("P12").Value = 0.65

$ws.Range("D13").Value = 0.86
$ws.Range("G13").Value = 0.8
$ws.Range("J13").Value = 0.86
$ws.Range("M13").Value = 0.86
$ws.Range("P13").Value = 0.86

$ws.Range("D14").Value = 0.86
$ws.Range("G14").Value = 0.61
$ws.Range("J14").Value = 0.62
$ws.Range("M14").Value = 0.72
$ws.Range("P14").Value = 0.8

$ws.Range("D15").Value = 0.87
$ws.Range("G15").Value = 0.76
$ws.Range("J15").Value = 0.8
$ws.Range("M15").Value = 0.83
$ws.Range("P15").Value = 0.85

